$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the original A1:E4 block (numbers + letters) to G1:K4 before
# making the in-place edits below -- this mirrors the "copy columns"
# feature the commit message refers to.
$ws.Range("A1:E4").Copy()
$ws.Range("G1").PasteSpecial()
$excel.CutCopyMode = 0

# In-place edits on the original block.
$ws.Range("A2").Value = "a<"
$ws.Range("B3").Value = "bb>"
$ws.Range("D1").Value = "?4"

# Match the saved selection state.
$ws.Range("I9").Select() | Out-Null
